$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-6 (years 2014-2018 IFRS columns shifted)
# Row 2
$ws.Range("D2").Value = 22134
$ws.Range("E2").Value = 638
$ws.Range("F2").Value = 638
$ws.Range("G2").Value = 463
$ws.Range("H2").Value = 210
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 185
$ws.Range("K2").Value = 27959
$ws.Range("L2").Value = 11416
$ws.Range("M2").Value = 16544
$ws.Range("N2").Value = 13151
$ws.Range("O2").Value = 3393
$ws.Range("P2").Value = 443
$ws.Range("Q2").Value = 1281
$ws.Range("R2").Value = -1856
$ws.Range("S2").Value = 1023
$ws.Range("T2").Value = 644
$ws.Range("U2").Value = 637
$ws.Range("V2").Value = 6924
$ws.Range("W2").Value = 2.88
$ws.Range("X2").Value = 0.95
$ws.Range("Y2").Value = 0.19
$ws.Range("Z2").Value = 0.78
$ws.Range("AA2").Value = 69
$ws.Range("AB2").Value = 2792.76
$ws.Range("AC2").Value = 296
$ws.Range("AD2").Value = 303.03
$ws.Range("AE2").Value = 170099
$ws.Range("AF2").Value = 0.53
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 1.67
$ws.Range("AI2").Value = 461.19
$ws.Range("AJ2").Value = 8564271

# Row 3
$ws.Range("D3").Value = 21696
$ws.Range("E3").Value = 1067
$ws.Range("F3").Value = 1067
$ws.Range("G3").Value = 1049
$ws.Range("H3").Value = 780
$ws.Range("I3").Value = 465
$ws.Range("J3").Value = 315
$ws.Range("K3").Value = 30714
$ws.Range("L3").Value = 12751
$ws.Range("M3").Value = 17964
$ws.Range("N3").Value = 13421
$ws.Range("O3").Value = 4543
$ws.Range("P3").Value = 443
$ws.Range("Q3").Value = 899
$ws.Range("R3").Value = -526
$ws.Range("S3").Value = -852
$ws.Range("T3").Value = 960
$ws.Range("U3").Value = -62
$ws.Range("V3").Value = 8096
$ws.Range("W3").Value = 4.92
$ws.Range("X3").Value = 3.6
$ws.Range("Y3").Value = 3.5
$ws.Range("Z3").Value = 2.66
$ws.Range("AA3").Value = 70.98
$ws.Range("AB3").Value = 2868.04
$ws.Range("AC3").Value = 5242
$ws.Range("AD3").Value = 30.33
$ws.Range("AE3").Value = 173590
$ws.Range("AF3").Value = 0.92
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 0.94
$ws.Range("AI3").Value = 24.98
$ws.Range("AJ3").Value = 8564271

# Row 4
$ws.Range("D4").Value = 23115
$ws.Range("E4").Value = 1410
$ws.Range("F4").Value = 1410
$ws.Range("G4").Value = 1027
$ws.Range("H4").Value = 651
$ws.Range("I4").Value = 226
$ws.Range("J4").Value = 425
$ws.Range("K4").Value = 31694
$ws.Range("L4").Value = 13188
$ws.Range("M4").Value = 18506
$ws.Range("N4").Value = 13526
$ws.Range("O4").Value = 4980
$ws.Range("P4").Value = 443
$ws.Range("Q4").Value = 2171
$ws.Range("R4").Value = -680
$ws.Range("S4").Value = -624
$ws.Range("T4").Value = 1303
$ws.Range("U4").Value = 867
$ws.Range("V4").Value = 7833
$ws.Range("W4").Value = 6.1
$ws.Range("X4").Value = 2.82
$ws.Range("Y4").Value = 1.68
$ws.Range("Z4").Value = 2.09
$ws.Range("AA4").Value = 71.26
$ws.Range("AB4").Value = 2881.71
$ws.Range("AC4").Value = 2550
$ws.Range("AD4").Value = 46.87
$ws.Range("AE4").Value = 174957
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 1.67
$ws.Range("AI4").Value = 68.45
$ws.Range("AJ4").Value = 8564271

# Row 5
$ws.Range("D5").Value = 23963
$ws.Range("E5").Value = 988
$ws.Range("F5").Value = 988
$ws.Range("G5").Value = 894
$ws.Range("H5").Value = 672
$ws.Range("I5").Value = 377
$ws.Range("J5").Value = 295
$ws.Range("K5").Value = 32760
$ws.Range("L5").Value = 13362
$ws.Range("M5").Value = 19398
$ws.Range("N5").Value = 13735
$ws.Range("O5").Value = 5565
$ws.Range("P5").Value = 443
$ws.Range("Q5").Value = -219
$ws.Range("R5").Value = -2422
$ws.Range("S5").Value = 1258
$ws.Range("T5").Value = 871
$ws.Range("U5").Value = -1090
$ws.Range("V5").Value = 8210
$ws.Range("W5").Value = 4.12
$ws.Range("X5").Value = 2.8
$ws.Range("Y5").Value = 2.76
$ws.Range("Z5").Value = 2.08
$ws.Range("AA5").Value = 68.88
$ws.Range("AB5").Value = 2928.91
$ws.Range("AC5").Value = 4248
$ws.Range("AD5").Value = 28.72
$ws.Range("AE5").Value = 177535
$ws.Range("AF5").Value = 0.69
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.64
$ws.Range("AI5").Value = 41.08
$ws.Range("AJ5").Value = 8564271

# Row 6
$ws.Range("D6").Value = 25635
$ws.Range("E6").Value = 1385
$ws.Range("F6").Value = 1385
$ws.Range("G6").Value = 1293
$ws.Range("H6").Value = 1025
$ws.Range("I6").Value = 700
$ws.Range("K6").Value = 33138
$ws.Range("L6").Value = 13119
$ws.Range("M6").Value = 20019
$ws.Range("N6").Value = 14116
$ws.Range("P6").Value = 443
$ws.Range("Q6").Value = 875
$ws.Range("R6").Value = -624
$ws.Range("S6").Value = -202
$ws.Range("T6").Value = 1583
$ws.Range("U6").Value = -708
$ws.Range("V6").Value = 7786
$ws.Range("W6").Value = 5.4
$ws.Range("X6").Value = 4
$ws.Range("Y6").Value = 5.03
$ws.Range("Z6").Value = 3.11
$ws.Range("AA6").Value = 65.53
$ws.Range("AB6").Value = 3047.14
$ws.Range("AC6").Value = 7890
$ws.Range("AD6").Value = 9.95
$ws.Range("AE6").Value = 182460
$ws.Range("AF6").Value = 0.43
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 2.55
$ws.Range("AI6").Value = 22.13
$ws.Range("AJ6").Value = 8564271

# Rows 7-9: clear all data columns (D:AJ), keep A:C labels intact
$ws.Range("D7:AJ9").ClearContents()
